$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose values are numeric-looking or percent-looking text (e.g. "304.87", "-0.72%").
# The source workbook stores these as literal text, so we force the Text number format
# before assigning, then restore the default "Normal" style so no stray formatting is left
# behind (matches the un-formatted cells elsewhere in the sheet).
$textCells = @{
    'D2' = '304.87'
    'E2' = '-0.72%'
    'D3' = '35.96'
    'E3' = '-1.50%'
    'D4' = '5.017'
    'E4' = '-1.99%'
    'D5' = '0.08065'
    'E5' = '-0.93%'
    'D6' = '1.883'
    'E6' = '-3.54%'
    'D7' = '7.835'
    'E7' = '0.76%'
    'D8' = '0.9308'
    'E8' = '-0.70%'
    'D9' = '0.1287'
    'E9' = '-12.90%'
    'D10' = '0.1905'
    'E10' = '-1.26%'
    'D11' = '0.09212'
    'D12' = '0.03512'
    'E12' = '-0.82%'
    'D13' = '0.09893'
    'E13' = '0.61%'
    'D14' = '0.001419'
    'E14' = '-1.09%'
    'D15' = '0.006359'
    'E15' = '9.27%'
    'E16' = '3.54%'
    'D17' = '4.150'
    'E17' = '0.43%'
    'D18' = '3.166'
    'E18' = '4.84%'
    'D19' = '0.3451'
    'E19' = '0.75%'
    'D20' = '0.1345'
    'E20' = '2.23%'
    'D21' = '5.223'
    'E21' = '4.77%'
    'E22' = '1.61%'
    'D23' = '0.04413'
    'E23' = '-2.24%'
    'E24' = '2.08%'
    'D25' = '0.004712'
    'E25' = '-3.91%'
    'D26' = '0.0001301'
    'E26' = '4.82%'
    'E27' = '-29.37%'
    'D39' = '0.01949'
    'E39' = '-2.01%'
    'D40' = '0.05153'
    'E40' = '5.86%'
    'D41' = '0.007537'
    'E41' = '-0.45%'
    'D42' = '0.01015'
    'E42' = '-7.20%'
    'D43' = '0.1372'
    'E43' = '-0.37%'
    'D44' = '0.002172'
    'E44' = '4.25%'
    'E45' = '11.32%'
    'D46' = '0.00006351'
    'E46' = '-0.39%'
    'E47' = '0.29%'
    'D48' = '63.57'
    'E48' = '-1.69%'
    'D49' = '0.001660'
    'E49' = '39.77%'
    'D50' = '0.00002102'
    'E50' = '0.29%'
    'D51' = '0.0002002'
    'E51' = '0.29%'
}

foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$addr]
    $cell.Style = "Normal"
}

# Plain text cells (coin names / hyperlink URLs) - no numeric auto-detection risk.
$textValues = @{
    'B20' = 'ProBitToken'
    'C20' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'B21' = 'MCDex'
    'C21' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
}

foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}
